$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: G3, H3 -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4 -> 2, E4 -> 1, F4 -> 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

# Row 5: D5 -> 1, E5 -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: D6 -> 1, E6 -> 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: D7 -> 1, E7 -> 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

# Row 8: H8 -> 1
$ws.Range("H8").Value = 1

# Row 9: D9 -> 1, E9 -> 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: D10 -> 1, E10 -> 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: H11 -> 1
$ws.Range("H11").Value = 1

# Row 12: D12 -> 1, E12 -> 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: D13 -> 2, E13 -> 1, F13 -> 1
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1

# Row 14: D14 -> 1, E14 -> 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1

# Row 15: H15 -> 1
$ws.Range("H15").Value = 1

# Row 16: H16 -> 1
$ws.Range("H16").Value = 1

# Row 17: H17 -> 1
$ws.Range("H17").Value = 1

# Row 18: D18 -> 1, E18 -> 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1
